$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Legend / text updates ---
$ws.Range("A18").Value = "N/A = test not relevant"
$ws.Range("A10").Value = "dbRDA*"
$ws.Range("A20").Value = "* Note: dbRDA is an ordination technique and does not have p-values. It produces ordination plots in accordance with the nature of the tests (i.e. positive/negative controls) and appears to produce meaningful results based on previous analysis using principal coordinates analysis. The percent variability that is explained by each dbRDA model's constraining variable is included in the table instead of p-values."

# --- dbRDA row (row 10): replace "TBD" placeholders with percent-variability values ---
$ws.Range("B10").Value = 0.2786
$ws.Range("C10").Value = 0.032887
$ws.Range("D10").Value = 0.006929
$ws.Range("E10").Value = 0.33187
$ws.Range("F10").Value = 0.0888133
$ws.Range("I10").Value = 0.2328
$ws.Range("J10").Value = 0.00957966

$ws.Range("B10:F10").NumberFormat = "0.00%"
$ws.Range("I10:J10").NumberFormat = "0.00%"

# --- Selection moved to P10 ---
$ws.Range("P10").Select()
